# Actualización automática 2025-10-17 08:30:08
#
# Updates sales figures for GUERRERO FAREZ FABIAN MAURICIO across the three
# worksheets (VENTAS POR GRUPO, VENTA MENSUAL, CUMPLIMIENTO MENSUAL) and
# refreshes the dependent subtotal / percentage cells that were stored as
# static values rather than live formulas.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: VENTAS POR GRUPO
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

$ws1.Range("M4").Value = 3670.77

$ws1.Range("D28").Value = 933.12
$ws1.Range("K28").Value = 1507.32
$ws1.Range("L28").Value = 760.3200000000001
$ws1.Range("M28").Value = 5819.2

$ws1.Range("M35").Value = 7699.44

$ws1.Range("M48").Value = 189.19

# Row 55 holds "<count> de 53" summaries of how many clients posted a
# positive amount in each product column. K, L and M each gained one more
# qualifying client because of the changes above.
$ws1.Range("K55").Value = "5 de 53"
$ws1.Range("L55").Value = "3 de 53"
$ws1.Range("M55").Value = "10 de 53"

# ---------------------------------------------------------------------
# Sheet 2: VENTA MENSUAL
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

$ws2.Range("F4").Value = 5290.76
$ws2.Range("F28").Value = 9019.959999999999
$ws2.Range("F35").Value = 9686.73
$ws2.Range("F48").Value = 189.19

# Row 59 is the column total; recompute it with the updated figures.
$ws2.Range("F59").Value = 41834.6

# ---------------------------------------------------------------------
# Sheet 3: CUMPLIMIENTO MENSUAL
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Column D got a little wider to fit the bigger numbers (the COM layer adds
# a constant ~0.8333 padding to whatever ColumnWidth is assigned, so the
# input is pre-compensated to land exactly on width="14" in the XML).
$ws3.Columns.Item(4).ColumnWidth = 13.166666666666666

# Row 3: 240X80 PORCELANATO
$ws3.Range("D3").Value = 5123.52
$ws3.Range("E3").Value = 12545.6270988183
$ws3.Range("F3").Value = 0.2899698537425532

# Row 10: PANELES DECORATIVOS
$ws3.Range("D10").Value = 2476.73
$ws3.Range("E10").Value = 1404.34983534392
$ws3.Range("F10").Value = 0.6381548705710987

# Row 11: PIEDRA SINTERIZADA
$ws3.Range("D11").Value = 5122.07
$ws3.Range("E11").Value = 6708.93
$ws3.Range("F11").Value = 0.432936353647198

# Row 12: PORCELANATO
$ws3.Range("D12").Value = 20562.57
$ws3.Range("E12").Value = 32100.55
$ws3.Range("F12").Value = 0.3904548382245488

# Row 14: TOTAL
$ws3.Range("D14").Value = 39813.91
$ws3.Range("E14").Value = 59202.59661190613
$ws3.Range("F14").Value = 0.4020936646053379
